# Code-review fix: the inventory end-date value had been typed into the
# wrong column (E2, under "tool_inventory_start_date") and belongs in F2
# (under "tool_inventory_end_date"). Move it there, keeping its value and
# number-format/style intact, then tidy up the sheet the way Excel would
# after such an edit (refresh the selection to the cell that now holds the
# value, and auto-fit the columns to their content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cut E2 (value + style) and paste it into F2; this removes the E2 cell
# entirely (no leftover empty/styled cell) and creates F2 with the same
# value/style that E2 had.
$ws.Range("E2").Cut($ws.Range("F2"))
$ws.Range("E2").Clear()

# After the move, Excel leaves the pasted-into cell selected.
$ws.Range("F2").Select()

# Auto-fit all the data columns to their (now updated) content.
$ws.Columns("A:F").AutoFit()
